# template_po_lien_cong_ty_ver2.0.xlsx -- "update fix code t11"
#
# The sheet "PO - HH" gets a new column inserted right before the existing
# column O ("Chi tiết đơn hàng / Hàng tặng"). The new column O carries a new
# header, "Chi tiết đơn hàng / Mô tả", and everything that used to live in
# O..AF shifts one column to the right (now O..AG). The cell note that was
# anchored on the old O1 has to move along with it, landing on the new P1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PO - HH")
$ws.Activate()

# Insert a new blank column before column O (column 15); this shifts every
# column from O onward (through AF) one position to the right, turning the
# old dimension A1:AF10 into A1:AG10 and keeping all existing values/styles
# intact on their new letters.
$ws.Columns.Item(15).Insert()

# The newly inserted column inherits its neighbour's formatting; just add
# the new header text for it.
$ws.Range("O1").Value = "Chi tiết đơn hàng / Mô tả"

# The note that used to sit on O1 ("Chi tiết đơn hàng / Hàng tặng" header)
# did not move with the insert, so relocate it by hand onto the cell that
# now holds that header, P1.
$oldNote = $ws.Range("O1").Comment
if ($oldNote -ne $null) {
    $noteText = $oldNote.Text()
    $oldNote.Delete()
    $ws.Range("P1").AddComment($noteText)
}
